# Kanban board update: move "Tarea 1" (Carlos Megias) from "In progress" (C6)
# to "Review" (D6), and leave the selection on the new location (D6).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$source = $ws.Range("C6")
$target = $ws.Range("D6")

# Move the task text to the Review column, keep each cell's own formatting.
$target.Value = $source.Value2
$source.ClearContents()

# Reflect the move in the current selection, as the author left it.
$target.Select()
